$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grids")
$ws.Range("G5:G596").Value = "lo"
